$d = $word.ActiveDocument

# --- Step 1: Fix paragraph 1's pPr (remove the eastAsia rFonts hint from the paragraph mark) ---
# Insert a throwaway paragraph (with the desired, hint-free pPr) right after paragraph 1's
# text, then delete the pilcrow that separates them. Deleting a paragraph mark merges the
# two paragraphs and the *surviving* mark (the later one) supplies the merged pPr, so
# paragraph 1 ends up with a clean pPr (spacing only, no rFonts hint).
$p1 = $d.Paragraphs.Item(1)
$insertPos = $p1.Range.End - 1
$ins1 = $d.Range($insertPos, $insertPos)
$xml1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:spacing w:line="220" w:lineRule="atLeast"/></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$ins1.InsertXML($xml1)
$p1b = $d.Paragraphs.Item(1)
$mark1 = $d.Range($p1b.Range.End - 1, $p1b.Range.End)
$mark1.Delete()

# --- Step 2: Extend paragraph 2 and append the five new paragraphs ---
# Insert the continuation text + new paragraphs as a new paragraph right after paragraph 2's
# existing text, then delete the pilcrow that used to end paragraph 2 so the continuation
# runs merge back into paragraph 2 (adopting the desired pPr carried on the inserted XML).
$p2 = $d.Paragraphs.Item(2)
$insertPos2 = $p2.Range.End - 1
$ins2 = $d.Range($insertPos2, $insertPos2)
$xml2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:spacing w:line="220" w:lineRule="atLeast"/><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> the requirments for job have changed </w:t></w:r><w:r><w:t>gradually</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>.Previously,people can finish their job as long as they master a skill.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="220" w:lineRule="atLeast"/><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>But now,situations have changed.Computer skill,language ability and knowledge about law are needed for an employee,which helps them go further in their career.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="220" w:lineRule="atLeast"/><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>Besides,when we attend the training class,we can meet many people in other fields.It</w:t></w:r><w:r><w:t>’</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>s a good chance to know something new about other fileds conveniently.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="220" w:lineRule="atLeast"/><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>And we can get more information about their fields that would be very helpful to our career.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="220" w:lineRule="atLeast"/><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>Finally,all of us should hold the idea that it</w:t></w:r><w:r><w:t>’</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>s never too late to lear.Learing is the only way to keep pace with the society and attending training clases may be the most effective way.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="220" w:lineRule="atLeast"/></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">From the </w:t></w:r><w:r><w:t>discussion</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> above,we  can conclude that it</w:t></w:r><w:r><w:t>’</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>s necessary and worthwhile for us to attend training classes because of its great importance.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$ins2.InsertXML($xml2)
$p2b = $d.Paragraphs.Item(2)
$mark2 = $d.Range($p2b.Range.End - 1, $p2b.Range.End)
$mark2.Delete()

Write-Host "Final paragraph count:" $d.Paragraphs.Count
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    Write-Host ("Para " + $i + ": " + $d.Paragraphs.Item($i).Range.Text)
}
